$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'27.100.89"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.17%  '

# Row 3
$ws.Range('D3').Value = "'1.818.49"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.38%  '

# Row 4
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '

# Row 5
$ws.Range('D5').Value = "'312.20"
$ws.Range('D5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = "'0.9995"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.15%  '

# Row 7
$ws.Range('D7').Value = "'0.4453"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.11%  '

# Row 8
$ws.Range('D8').Value = "'0.3738"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.85%  '

# Row 9
$ws.Range('D9').Value = "'0.07473"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.46%  '

# Row 10
$ws.Range('D10').Value = "'0.8717"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.10%  '

# Row 11
$ws.Range('D11').Value = "'20.87"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.92%  '

# Row 12
$ws.Range('D12').Value = "'1.810.93"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.81%  '

# Row 13
$ws.Range('D13').Value = "'6.735"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.10%  '

# Row 14
$ws.Range('D14').Value = "'94.35"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.37%  '

# Row 15
$ws.Range('D15').Value = "'5.343"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.94%  '

# Row 16
$ws.Range('D16').Value = "'0.07093"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.84%  '

# Row 17
$ws.Range('D17').Value = "'0.9999"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.31%  '

# Row 18
$ws.Range('D18').Value = "'0.000008760"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.06%  '

# Row 19
$ws.Range('D19').Value = "'0.9992"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.19%  '

# Row 20
$ws.Range('E20').Value = '  +0.91%  '

# Row 21
$ws.Range('D21').Value = "'27.138.16"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.12%  '

# Row 22
$ws.Range('D22').Value = "'5.228"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.87%  '

# Row 24
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = "'1.988"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.60%  '

# Row 25
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = "'2.426"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.33%  '

# Row 26
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'151.41"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.27%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'18.51"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.68%  '

# Row 28
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = "'5.310"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.11%  '

# Row 29
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'117.93"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.01%  '

# Row 30
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = "'0.08814"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.97%  '

# Row 31
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'0.7676"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.08%  '

# Row 32
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = "'1.176"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.23%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'4.557"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.87%  '

# Row 34
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = "'2.887"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.48%  '

# Row 35
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = "'0.9990"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.16%  '

# Row 36
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = "'1.104"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.08%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.01984"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.10%  '

# Row 38
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = "'7.439"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.36%  '

# Row 39
$ws.Range('D39').Value = "'0.05273"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.64%  '

# Row 40
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.5318"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.85%  '

# Row 41
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = "'2.855"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.69%  '

# Row 42
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.1715"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.64%  '

# Row 43
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'2.188"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.36%  '

# Row 44
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'8.722"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.94%  '

# Row 45
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.5049"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.66%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'10.59"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.71%  '

# Row 47
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = "'1.704"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.24%  '

# Row 48
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'105.52"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.27%  '

# Row 49
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'0.9988"
$ws.Range('D49').Style = 'Normal'

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.06373"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.76%  '

# Row 51
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = "'0.9307"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.09%  '
